$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 23:35"

# 2. Update Estados Unidos row (row 4) stats
$ws.Range("B4").Value = 1405915
$ws.Range("C4").Value = 20081
$ws.Range("D4").Value = 276440
$ws.Range("E4").Value = 1046195
$ws.Range("G4").Value = 1485
$ws.Range("H4").Value = 83280

# 3. Gabon's updated numbers push it above Niger and Costa Rica in the
#    ranking (sorted descending by total cases, column B). Re-assign the
#    country names and stats for rows 107-109 so the table stays sorted.
$ws.Range("A107").Value = "Gabon"
$ws.Range("B107").Value = 863
$ws.Range("C107").Value = 61
$ws.Range("D107").Value = 137
$ws.Range("E107").Value = 717
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 9

$ws.Range("A108").Value = "Niger"
$ws.Range("B108").Value = 854
$ws.Range("C108").Value = 22
$ws.Range("D108").Value = 648
$ws.Range("E108").Value = 159
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 47

$ws.Range("A109").Value = "Costa Rica"
$ws.Range("B109").Value = 804
$ws.Range("C109").Value = 3
$ws.Range("D109").Value = 520
$ws.Range("E109").Value = 277
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 7

# 4. Row 154 updates
$ws.Range("B154").Value = 148
$ws.Range("C154").Value = 1
$ws.Range("D154").Value = 144

# 5. Row 160 updates
$ws.Range("B160").Value = 121
$ws.Range("C160").Value = 2
$ws.Range("E160").Value = 47
$ws.Range("F160").Value = 2
